# Air_WorldSpan.xlsx update:
#  - Supplier column (N) changes from "Worldspan" to "WSpan" on rows that
#    reflect the supplier-name change in the release.
#      Sheet "Air_WorldSpan_OneWay"      -> N3, N4, N5 (N2 left as-is)
#      Sheet "Air_WorldSpan_RoundTrip"   -> N2, N3, N4, N5
#  - Active sheet / tab selection flips from the RoundTrip sheet to the
#    OneWay sheet (also updates each sheet's remembered scroll / selection).

$wb = $excel.ActiveWorkbook

$wsOneWay   = $wb.Worksheets.Item("Air_WorldSpan_OneWay")
$wsRoundTrip = $wb.Worksheets.Item("Air_WorldSpan_RoundTrip")

# --- Update Supplier values -------------------------------------------------
$wsOneWay.Range("N3").Value = "WSpan"
$wsOneWay.Range("N4").Value = "WSpan"
$wsOneWay.Range("N5").Value = "WSpan"

$wsRoundTrip.Range("N2").Value = "WSpan"
$wsRoundTrip.Range("N3").Value = "WSpan"
$wsRoundTrip.Range("N4").Value = "WSpan"
$wsRoundTrip.Range("N5").Value = "WSpan"

# --- Update view / selection state ------------------------------------------
# RoundTrip: scroll/selection moves to N5 (tabSelected is dropped here since
# OneWay becomes the selected tab instead).
$wsRoundTrip.Activate()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$wsRoundTrip.Range("N5").Select()

# OneWay: becomes the selected/active tab, scrolled/selected at N4.
$wsOneWay.Activate()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 3
$wsOneWay.Range("N4").Select()
